# Applies the scheduled-runner price/profit recalculation update across
# the ALC, ARM, CRP, CUL, GSM and WVR sheets (columns H-N: average price,
# average price NQ/HQ, leve price NQ/HQ, leve profit NQ/HQ).

$wb = $excel.ActiveWorkbook

function Set-Row {
    param(
        [string]$SheetName,
        [int]$Row,
        [hashtable]$Values   # column letter -> value ($null to clear)
    )
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($col in $Values.Keys) {
        $val = $Values[$col]
        $rng = $ws.Range("$col$Row")
        if ($null -eq $val) {
            $rng.ClearContents()
        } else {
            $rng.Value = $val
        }
    }
}

# ---- ALC ----
Set-Row "ALC" 141 @{
    H = 1215.3
    I = 935.4595
    J = 4666.6665
    K = 2806.3785
    L = 13999.9995
    M = 2373.6215
    N = -24359.9995
}

# ---- ARM ----
Set-Row "ARM" 88 @{
    H = 5283.3335
    I = 2366.6667
    J = 6255.5557
    K = 2366.6667
    L = 6255.5557
    M = -1960.6667
    N = -7067.5557
}

Set-Row "ARM" 91 @{
    H = 5283.3335
    I = 2366.6667
    J = 6255.5557
    K = 2366.6667
    L = 6255.5557
    M = -962.6667000000002
    N = -9063.555700000001
}

# ---- CRP ----
Set-Row "CRP" 62 @{
    H = 2587.2
    I = 2249.2856
    K = 2249.2856
    M = -1625.2856
}

Set-Row "CRP" 65 @{
    H = 2587.2
    I = 2249.2856
    K = 11246.428
    M = -8126.428
}

# ---- CUL ----
Set-Row "CUL" 5 @{
    H = 2828090
    J = 1334553.2
    L = 4003659.6
    N = -4003883.6
}

Set-Row "CUL" 12 @{
    H = 70.22727
    J = 74.73333
    L = 224.19999
    N = -570.19999
}

Set-Row "CUL" 22 @{
    H = 1750
    I = 1050
    J = 1983.3334
    K = 3150
    L = 5950.0002
    M = -2981
    N = -6288.0002
}

Set-Row "CUL" 27 @{
    H = 1750
    I = 1050
    J = 1983.3334
    K = 3150
    L = 5950.0002
    M = -3048
    N = -6154.0002
}

Set-Row "CUL" 40 @{
    H = 6763.5293
    I = 7877.143
    J = 1566.6666
    K = 31508.572
    L = 6266.6664
    M = -31439.572
    N = -6404.6664
}

Set-Row "CUL" 46 @{
    H = 2401.8572
    I = 1515
    J = 2756.6
    K = 4545
    L = 8269.799999999999
    M = -4454
    N = -8451.799999999999
}

Set-Row "CUL" 58 @{
    H = 2542.5
    I = 90
    J = 2892.8572
    K = 270
    L = 8678.571599999999
    M = -142
    N = -8934.571599999999
}

Set-Row "CUL" 64 @{
    H = 1762.6
    I = 375.14285
    J = 5000
    K = 1125.42855
    L = 15000
    M = -855.4285500000001
    N = -15540
}

Set-Row "CUL" 67 @{
    H = 1762.6
    I = 375.14285
    J = 5000
    K = 1125.42855
    L = 15000
    M = -189.4285500000001
    N = -16872
}

Set-Row "CUL" 70 @{
    H = 2061.25
    I = 1248.3334
    J = 4500
    K = 3745.0002
    L = 13500
    M = -3430.0002
    N = -14130
}

Set-Row "CUL" 73 @{
    H = 2061.25
    I = 1248.3334
    J = 4500
    K = 3745.0002
    L = 13500
    M = -2653.0002
    N = -15684
}

Set-Row "CUL" 76 @{
    H = 0
    I = 0
    J = 0
    K = 0
    L = 0
    M = $null
    N = $null
}

Set-Row "CUL" 79 @{
    H = 0
    I = 0
    J = 0
    K = 0
    L = 0
    M = $null
    N = $null
}

Set-Row "CUL" 82 @{
    H = 26644.445
    I = 700
    J = 29887.5
    K = 2100
    L = 89662.5
    M = -1694
    N = -90474.5
}

Set-Row "CUL" 85 @{
    H = 26644.445
    I = 700
    J = 29887.5
    K = 2100
    L = 89662.5
    M = -696
    N = -92470.5
}

Set-Row "CUL" 88 @{
    H = 13000
    J = 13000
    L = 39000
    N = -39856
}

Set-Row "CUL" 91 @{
    H = 13000
    J = 13000
    L = 39000
    N = -41964
}

Set-Row "CUL" 122 @{
    H = 1721.6842
    I = 300.4
    J = 2229.2856
    K = 2703.6
    L = 20063.5704
    M = -253.5999999999999
    N = -24963.5704
}

Set-Row "CUL" 135 @{
    H = 2828090
    J = 1334553.2
    L = 12010978.8
    N = -12016048.8
}

# ---- GSM ----
Set-Row "GSM" 70 @{
    H = 3791942.2
    I = 2061849.8
    J = 6539736
    K = 2061849.8
    L = 6539736
    M = -2061579.8
    N = -6540276
}

Set-Row "GSM" 73 @{
    H = 3791942.2
    I = 2061849.8
    J = 6539736
    K = 2061849.8
    L = 6539736
    M = -2060913.8
    N = -6541608
}

# ---- WVR ----
Set-Row "WVR" 2 @{
    H = 9650
    I = 9650
    K = 9650
    M = -9538
}
